$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Seed the new "Pass" shared string first (so it lands at the same shared
#    string index the canonical file uses), then update the row 7 wording
#    from the old "i_countenb" phrasing to the new "i_rtcclk" phrasing.
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = "Pass"

$ws.Range("D7").Value = "i_rtcclk' of every counter is connected to the flag output 'o_rolloverflag' of the previous counter"
$ws.Range("E7").Value = "observe 'i_rtcclk' and compare to the respective 'o_rolloverflag'"
$ws.Range("G7").Value = "i_rtcclk' of every subcounter = previous 'o_rolloverflag' except first one 'i_rtcclk' from trigger detection"

# Re-applying the value above resets the quote-prefix flag that D7 / G7 used
# to carry, so restore it by pulling the formatting (only) from neighbouring
# cells that already carry the quote-prefix + wrap style.
$ws.Range("D8").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Populate the new "Actual Results" column (H) as a copy of the
#    "Expected results / Outputs" column (G), values + formatting alike.
# ---------------------------------------------------------------------------
$ws.Range("G2:G16").Copy($ws.Range("H2:H16"))

# Two rows (6 and 14) pick up word-wrap in column H that column G itself
# never had.
$ws.Range("H6").WrapText = $true
$ws.Range("H14").WrapText = $true

# ---------------------------------------------------------------------------
# 3) Populate the new "Pass/Fail" column (I) with "Pass" for every test row.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 9).Value = "Pass"
    $ws.Cells.Item($r, 9).WrapText = $true
}

# ---------------------------------------------------------------------------
# 4) Row-height touch ups that Excel re-flowed once the new text landed.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 72
$ws.Rows.Item(16).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 5) Selection / scroll position as last left by the editor.
# ---------------------------------------------------------------------------
$ws.Range("B15").Select()
$excel.ActiveWindow.ScrollRow = 1
